# Api Endpoints.xlsx - add new routes for #17:
#   POST /, POST /users, GET /users, GET /users/{id}, and PUT /users/{id}, POST /login.
#
# In the existing sheet this shows up as:
#   - the GET /users row (C46) becomes GET /users?type={type-id}
#   - the POST /users/{user-id}/login row (C49) becomes POST /login
#   - the /artists group header rows (24-27) and the /users group header rows (45-49)
#     get the green "group" marker in column A (same style already used by the
#     groups above them, e.g. A3, A5:A7, A9:A12, ...)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routes")

# --- text updates -------------------------------------------------------
$ws.Range("C46").Value = "/users?type={type-id}"
$ws.Range("C49").Value = "/login"

# --- green group markers in column A ------------------------------------
# 0 + 176*256 + 80*65536 = RGB(0, 176, 80), the same green already used by
# the other endpoint-group markers (e.g. A3, A5:A7, A9:A12, ...).
$markerColor = 5287936
$ws.Range("A24:A27").Interior.Color = $markerColor
$ws.Range("A45:A49").Interior.Color = $markerColor

# --- view state, matches where the author left the selection ------------
$ws.Activate()
$ws.Range("C23").Select()
